$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# ---- String constants ----
$s4 = @'
evaluator_partial_correctness
'@

$s5 = @'
 Given is the adjacency matrix for a weighted directed graph containing 16 nodes labelled A to P. The value corresponding to each row M and column N represents the cost of travelling between the two nodes, where 0 means no connection.   
what is the least cost path from node A to node P? Return the sequence of nodes in response.
   A B C D E F G H I J K L M N O P
 A 0 0 0 0 4 0 0 0 0 0 0 0 0 0 0 0
 B 2 0 1 0 0 0 0 0 0 0 0 0 0 0 0 0
 C 0 0 0 5 0 0 0 0 0 0 0 0 0 0 0 0
 D 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 E 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 F 0 5 0 0 4 0 0 0 0 2 0 0 0 0 0 0
 G 0 0 4 0 0 1 0 4 0 0 0 0 0 0 0 0
 H 0 0 0 1 0 0 0 0 0 0 0 0 0 0 0 0
 I 0 0 0 0 4 0 0 0 0 0 0 0 2 0 0 0
 J 0 0 0 0 0 0 0 0 5 0 0 0 0 0 0 0
 K 0 0 0 0 0 0 0 0 0 3 0 2 0 0 1 0
 L 0 0 0 0 0 0 0 3 0 0 0 0 0 0 0 0
 M 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 N 0 0 0 0 0 0 0 0 0 2 0 0 2 0 3 0
 O 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 P 0 0 0 0 0 0 0 0 0 0 0 3 0 0 2 0
    
'@

$s6 = @'
No possible path from A to P
'@

$s7 = @'
From the adjacency matrix, there seems to be a mistake as there is no direct or indirect link between node A and node P. Therefore, we cannot find a path from A to P.
'@

$s8 = @'
invalid input
'@

$s9 = @'
1/1
'@

$s10 = @'
 Given is the adjacency matrix for a weighted directed graph containing 25 nodes labelled A to Y. The value corresponding to each row M and column N represents the cost of travelling between the two nodes, where 0 means no connection.   
what is the least cost path from node A to node Y? Return the sequence of nodes in response.
   A B C D E F G H I J K L M N O P Q R S T U V W X Y
 A 0 5 0 0 0 2 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 B 0 0 3 0 0 0 1 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 C 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 D 0 0 2 0 0 0 0 0 2 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 E 0 0 0 4 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 F 0 0 0 0 0 0 5 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 G 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 H 0 0 0 0 0 0 4 0 2 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 I 0 0 0 0 0 0 0 0 0 1 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 J 0 0 0 0 5 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 K 0 0 0 0 0 4 0 0 0 0 0 0 0 0 0 5 0 0 0 0 0 0 0 0 0
 L 0 0 0 0 0 0 0 0 0 0 3 0 3 0 0 0 1 0 0 0 0 0 0 0 0
 M 0 0 0 0 0 0 0 1 0 0 0 0 0 5 0 0 0 0 0 0 0 0 0 0 0
 N 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 O 0 0 0 0 0 0 0 0 0 5 0 0 0 4 0 0 0 0 0 2 0 0 0 0 0
 P 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 5 0 0 0 0 0 0 0 0
 Q 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 R 0 0 0 0 0 0 0 0 0 0 0 0 5 0 0 0 4 0 3 0 0 0 0 0 0
 S 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 3 0 0 0 0 0
 T 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 U 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 5 0 0 0 0 0 0 0 0 0
 V 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 4 0 0 0 5 0 0 0 0
 W 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 4 0 0 0 3 0 1 0
 X 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 4
 Y 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 3 0 0 0 0 0
    
'@

$s11 = @'
No possible path from A to Y
'@

$s12 = @'
The adjacency matrix you provided suggests that there is no direct path from Node A to Node Y. Therefore, we have to consider other paths as well.
The following can be a possible path:
A -> B (cost 5), B -> C (cost 3), C -> D (cost 2), D -> I (cost 1), I -> Y (cost 4)
Total cost = 5+3+2+1+4 = 15
There could be other paths as well, but without using an algorithm or tracing all possible paths and their associated costs it may not be possible to definitively find the smallest cost path. I have not considered paths that have 0 cost connections as you mentioned 0 means no connection.
'@

$s13 = @'
0/1
'@

$s14 = @'
 Given is the adjacency matrix for a weighted directed graph containing 25 nodes labelled A to Y. The value corresponding to each row M and column N represents the cost of travelling between the two nodes, where 0 means no connection.   
what is the least cost path from node A to node Y? Return the sequence of nodes in response.
   A B C D E F G H I J K L M N O P Q R S T U V W X Y
 A 0 3 0 0 0 3 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 B 0 0 2 0 0 0 1 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 C 0 0 0 0 0 0 0 1 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 D 0 0 2 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 E 0 0 0 5 0 0 0 0 0 2 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 F 0 0 0 0 0 0 2 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 G 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 H 0 0 0 0 0 0 2 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 I 0 0 0 0 0 0 0 2 0 0 0 0 0 4 0 0 0 0 0 0 0 0 0 0 0
 J 0 0 0 0 0 0 0 0 1 0 0 0 0 0 1 0 0 0 0 0 0 0 0 0 0
 K 0 0 0 0 0 3 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 L 0 0 0 0 0 0 2 0 0 0 2 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 M 0 0 0 0 0 0 0 5 0 0 0 4 0 0 0 0 0 1 0 0 0 0 0 0 0
 N 0 0 0 0 0 0 0 0 0 0 0 0 4 0 0 0 0 0 0 0 0 0 0 0 0
 O 0 0 0 0 0 0 0 0 0 0 0 0 0 1 0 0 0 0 0 4 0 0 0 0 0
 P 0 0 0 0 0 0 0 0 0 0 1 0 0 0 0 0 0 0 0 0 4 0 0 0 0
 Q 0 0 0 0 0 0 0 0 0 0 0 4 0 0 0 5 0 0 0 0 0 1 0 0 0
 R 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 5 0 0 0 0 0 0 0 0
 S 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 4 0 3 0 0 0 3 0
 T 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 4
 U 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 2 0 0 0
 V 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 5 0 0
 W 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 1 0 0 0 0 0 0 0
 X 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 1 0 0
 Y 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 5 0
    
'@

$s15 = @'
From the given matrix, we can quickly identify the lowest-cost path from node A to node Y by following the lowest numbers in each row, moving from left to right.
A -> B (Cost: 3)
B -> G (Cost: 1)
G -> There are no paths leading from G to other nodes. This is a dead end.
Instead of going from B to G, we should go from B to C.
A -> B (Cost: 3)
B -> C (Cost: 2)
C -> H (Cost: 1)
H -> There are no paths leading from H to other nodes. Again, a dead end. 
Considering the next lower cost from B, 
A -> B (Cost: 3)
B -> There is no other path from B. This seems to be a wrong path to reach Y.
Let's try the next lower cost option A-> F.
A -> F (Cost: 3)
F -> I (Cost : 2)
I -> N (Cost : 4)
N -> Q (Cost : 4)
Q -> U (Cost : 1)
U -> Y (Cost : 2)
The total cost is 3+2+4+4+1+2 = 16 units.
So, the least cost path from node A to node Y is A -> F -> I -> N -> Q -> U -> Y.
'@

$s16 = @'
0/0
'@


# ---- Sheet1 (o_10): add new column E header (copy D1 formatting) + replace row2 content ----
$ws1.Range("D1:D1").Copy($ws1.Range("E1"))
$ws1.Range("E1").Value2 = $s4
$ws1.Range("A2").Value2 = $s5
$ws1.Range("B2").Value2 = $s6
$ws1.Range("C2").Value2 = $s7
$ws1.Range("D2").Value2 = $s8
$ws1.Range("E2").Value2 = $s9
$ws1.Rows.Item(2).AutoFit()

# ---- Create sheet2 (o_20) as a copy of sheet1 (keeps header formatting) ----
$ws1.Copy([System.Reflection.Missing]::Value, $ws1)
$ws2 = $wb.Worksheets.Item(2)
$ws2.Name = "o_20"
$ws2.Range("A2").Value2 = $s10
$ws2.Range("B2").Value2 = $s11
$ws2.Range("C2").Value2 = $s12
$ws2.Range("D2").Value2 = $s8
$ws2.Range("E2").Value2 = $s13
$ws2.Rows.Item(2).AutoFit()

# ---- Create sheet3 (o_20_jumbled) as a copy of sheet2, placed after it ----
$ws2.Copy([System.Reflection.Missing]::Value, $ws2)
$ws3 = $wb.Worksheets.Item(3)
$ws3.Name = "o_20_jumbled"
$ws3.Range("A2").Value2 = $s14
$ws3.Range("B2").Value2 = $s11
$ws3.Range("C2").Value2 = $s15
$ws3.Range("D2").Value2 = $s8
$ws3.Range("E2").Value2 = $s16
$ws3.Rows.Item(2).AutoFit()

# ---- Restore active/selected tab to sheet1 ----
$ws1.Activate()
